# Updated symbol list on Sun Jan 22 09:59:34 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the latest crypto snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values are written as Text (matching the
# original inline-string cells) so things like "302.04" and "-1.25%" are not
# auto-converted into numbers/percentages by Excel.
$updates = [ordered]@{
    "D2" = "302.04"
    "E2" = "-1.25%"
    "D3" = "37.27"
    "E3" = "6.69%"
    "D4" = "4.989"
    "E4" = "-4.03%"
    "D5" = "0.07810"
    "E5" = "-0.61%"
    "D6" = "2.189"
    "E6" = "-7.76%"
    "D7" = "8.033"
    "E7" = "-0.02%"
    "D8" = "4.043"
    "E8" = "2.27%"
    "D9" = "0.9153"
    "E9" = "-1.73%"
    "D10" = "0.09727"
    "E10" = "-3.97%"
    "D11" = "0.1884"
    "E11" = "2.36%"
    "D12" = "0.08678"
    "E12" = "0.20%"
    "D13" = "0.03544"
    "E13" = "6.40%"
    "D14" = "0.09971"
    "E14" = "0.79%"
    "D15" = "0.001486"
    "E15" = "0.39%"
    "D16" = "0.005648"
    "E16" = "0.55%"
    "D17" = "3.459"
    "E17" = "-0.80%"
    "D18" = "2.366"
    "E18" = "11.26%"
    "E19" = "2.18%"
    "D20" = "0.1276"
    "E20" = "-2.02%"
    "D21" = "4.778"
    "E21" = "10.72%"
    "D23" = "0.04637"
    "E23" = "1.46%"
    "D24" = "0.001231"
    "E24" = "1.25%"
    "D25" = "0.004791"
    "E25" = "7.70%"
    "E26" = "-7.31%"
    "E27" = "39.80%"
    "D39" = "0.01758"
    "E39" = "-1.66%"
    "D40" = "0.04736"
    "E40" = "-1.72%"
    "D41" = "0.008063"
    "E41" = "4.44%"
    "D42" = "0.1392"
    "E42" = "-1.42%"
    "D43" = "0.007695"
    "E43" = "8.96%"
    "D44" = "0.002136"
    "E44" = "-3.02%"
    "D45" = "0.01046"
    "E45" = "9.84%"
    "D46" = "0.00006058"
    "E46" = "2.10%"
    "E47" = "0.57%"
    "D48" = "7.911"
    "E48" = "189.52%"
    "D50" = "0.00002104"
    "E50" = "0.57%"
    "D51" = "0.0002004"
    "E51" = "0.57%"
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"      # force text so numeric-looking strings stay text
    $rng.Value = $updates[$cell]
    $rng.ClearFormats()          # drop the temporary "@" style override again
}
